$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural inserts (done first, while row numbers are still predictable) ---
# new blank row for AgPlenus / P23 (goes between existing P210 and P211 rows)
$ws.Rows("4:4").Insert()
# new blank row for LavieBio / P144 (goes between existing P145 and P143 rows)
$ws.Rows("7:7").Insert()

# --- fill in the new cell values ---
# LavieBio / P144 row
$ws.Range("B7").Value2 = "P144"
$ws.Range("A7").Value2 = "LavieBio"

# Upkeep / P997 row, appended after the last existing data row (row 40)
$ws.Range("B41").Value2 = "P997"
$ws.Range("A41").Value2 = "Upkeep"

# AgPlenus / P23 row
$ws.Range("B4").Value2 = "P23"
$ws.Range("A4").Value2 = "AgPlenus"

# match the saved selection state
$ws.Range("B5").Select()
